# Network and GSEA updating
$wb = $excel.ActiveWorkbook

$wsRelated = $wb.Worksheets.Item("Related")

# Swap the "Condition" (C) and "TimeCourse" (D) columns on the "Related" sheet,
# including the header row, for every populated row (1:49).
$colC = $wsRelated.Range("C1:C49").Value()
$colD = $wsRelated.Range("D1:D49").Value()
$wsRelated.Range("C1:C49").Value = $colD
$wsRelated.Range("D1:D49").Value = $colC

# Column widths on "Related" so the newly relabeled columns fit their content
# (matches the best-fit widths Excel computed for A:E after the edit).
$wsRelated.Columns.Item(1).ColumnWidth = 23.166666666666668
$wsRelated.Columns.Item(2).ColumnWidth = 9.5
$wsRelated.Columns.Item(3).ColumnWidth = 8.666666666666666
$wsRelated.Columns.Item(4).ColumnWidth = 10.666666666666666
$wsRelated.Columns.Item(5).ColumnWidth = 9.0

# Make "Related" the active sheet with its new selection.
$wsRelated.Activate()
$wsRelated.Range("G9").Select() | Out-Null
